$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.951.68"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "3.484.85"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.50"
$ws.Range("E5").Value = "  +6.32%  "
$ws.Range("E6").Value = "  -5.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  +3.65%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.01"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "4.043.12"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "3.487.60"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.28"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.06"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "65.930.97"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.25"
$ws.Range("E21").Value = "  -4.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("E22").Value = "  +9.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.40"
$ws.Range("E23").Value = "  +6.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.75"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.46"
$ws.Range("E25").Value = "  +10.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.05"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.22"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.33"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.69"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.76"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "595.30"
$ws.Range("E32").Value = "  -8.32%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.00"
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.153"
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0794"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.97"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.385"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.225.88"
$ws.Range("E41").Value = "  +6.28%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.55"
$ws.Range("E45").Value = "  -5.27%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0420"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.133"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  -4.89%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.65"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.01"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("E51").Value = "  -3.36%  "
